$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D7").Value = -7.269999999999999
$ws.Range("A10").Value = -20.926
$ws.Range("A12").Value = -21.694
$ws.Range("B13").Value = 6.606
$ws.Range("A18").Value = -21.694
$ws.Range("D20").Value = -8.222
